# "Basic the sistem is finished" -- replace the customer rows with the
# newly-collected leads data (only 2 data rows remain, a new "Observações"
# style column O is added, and a few fields are freeform text that must not
# be auto-converted to dates/numbers by Excel).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Eliseu Miguel ---
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Eliseu Miguel"
$ws.Range("C2").Value = "Itabuna"
$ws.Range("D2").Value = "Ba"
$ws.Range("E2").Value = "73 991920444"
$ws.Range("F2").Value = "'10/06/2002"
$ws.Range("G2").Value = "mguelmarinho@hotmail.com"
$ws.Range("H2").Value = "'3900"
$ws.Range("I2").Value = "Interessado"
$ws.Range("J2").Value = "Apartamento"
$ws.Range("K2").Value = "Imóvel Próprio"
$ws.Range("L2").Value = "36 meses"
$ws.Range("M2").Value = "R$250.000,00"
$ws.Range("N2").Value = "Financimaneto"
$ws.Range("O2").Value = "None"

# --- Row 3: Elias Gabriel Marinho de Oliveira ---
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Elias Gabriel Marinho de Oliveira"
$ws.Range("C3").Value = "Itabuna"
$ws.Range("D3").Value = "Ba"
$ws.Range("E3").Value = "73 991726631"
$ws.Range("F3").Value = "'29/07/1998"
$ws.Range("G3:O3").Value = ""

# --- Drop the old rows 4-6 (no longer part of the sheet) ---
$ws.Range("A4:O6").Delete()
